# Scheduled runner update: refresh market-price-derived profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on a handful of leve
# rows across several crafting-job sheets, per the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 313
$ws.Range("I2").Value = 194.5
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 194.5
$ws.Range("L2").Value = 550
$ws.Range("M2").Value = -81.5
$ws.Range("N2").Value = -776

$ws.Range("H21").Value = 12466.667
$ws.Range("J21").Value = 13960
$ws.Range("L21").Value = 13960
$ws.Range("N21").Value = -14896

$ws.Range("H23").Value = 12466.667
$ws.Range("J23").Value = 13960
$ws.Range("L23").Value = 13960
$ws.Range("N23").Value = -14428

$ws.Range("H38").Value = 2016485
$ws.Range("I38").Value = 2304254.2
$ws.Range("J38").Value = 2100
$ws.Range("K38").Value = 6912762.600000001
$ws.Range("L38").Value = 6300
$ws.Range("M38").Value = -6912390.600000001
$ws.Range("N38").Value = -7044

$ws.Range("H43").Value = 1899.6
$ws.Range("I43").Value = 6000.5
$ws.Range("J43").Value = 874.375
$ws.Range("K43").Value = 6000.5
$ws.Range("L43").Value = 874.375
$ws.Range("M43").Value = -5931.5
$ws.Range("N43").Value = -1012.375

$ws.Range("H58").Value = 745009.9
$ws.Range("I58").Value = 1116314.8
$ws.Range("J58").Value = 2400
$ws.Range("K58").Value = 3348944.4
$ws.Range("L58").Value = 7200
$ws.Range("M58").Value = -3348794.4
$ws.Range("N58").Value = -7500

$ws.Range("H61").Value = 2211095.8
$ws.Range("I61").Value = 2857261.2
$ws.Range("J61").Value = 57211.332
$ws.Range("K61").Value = 8571783.600000001
$ws.Range("L61").Value = 171633.996
$ws.Range("M61").Value = -8571611.600000001
$ws.Range("N61").Value = -171977.996

$ws.Range("H76").Value = 14532.667
$ws.Range("I76").Value = 3966
$ws.Range("J76").Value = 35666
$ws.Range("K76").Value = 3966
$ws.Range("L76").Value = 35666
$ws.Range("M76").Value = -3651
$ws.Range("N76").Value = -36296

$ws.Range("H79").Value = 14532.667
$ws.Range("I79").Value = 3966
$ws.Range("J79").Value = 35666
$ws.Range("K79").Value = 3966
$ws.Range("L79").Value = 35666
$ws.Range("M79").Value = -2874
$ws.Range("N79").Value = -37850

$ws.Range("H87").Value = 31270.2
$ws.Range("J87").Value = 31270.2
$ws.Range("L87").Value = 31270.2
$ws.Range("N87").Value = -33766.2

$ws.Range("H90").Value = 31270.2
$ws.Range("J90").Value = 31270.2
$ws.Range("L90").Value = 93810.60000000001
$ws.Range("N90").Value = -106290.6

$ws.Range("H107").Value = 551.6667
$ws.Range("I107").Value = 583.63635
$ws.Range("K107").Value = 583.63635
$ws.Range("M107").Value = 1336.36365

$ws.Range("H141").Value = 3978
$ws.Range("I141").Value = 3666.6667
$ws.Range("J141").Value = 4445
$ws.Range("K141").Value = 11000.0001
$ws.Range("L141").Value = 13335
$ws.Range("M141").Value = -5820.000100000001
$ws.Range("N141").Value = -23695

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2281.5283
$ws.Range("I132").Value = 2108.8333
$ws.Range("J132").Value = 2940.9092
$ws.Range("K132").Value = 6326.499899999999
$ws.Range("L132").Value = 8822.7276
$ws.Range("M132").Value = -3796.499899999999
$ws.Range("N132").Value = -13882.7276

$ws.Range("H135").Value = 42443.625
$ws.Range("J135").Value = 42443.625
$ws.Range("L135").Value = 42443.625
$ws.Range("N135").Value = -52583.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 113967.555
$ws.Range("I20").Value = 145744
$ws.Range("J20").Value = 2750
$ws.Range("K20").Value = 145744
$ws.Range("L20").Value = 2750
$ws.Range("M20").Value = -145497
$ws.Range("N20").Value = -3244

$ws.Range("H26").Value = 32666.666
$ws.Range("I26").Value = 10000
$ws.Range("J26").Value = 44000
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 44000
$ws.Range("M26").Value = -9708
$ws.Range("N26").Value = -44584

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 745.7059
$ws.Range("I114").Value = 311.2
$ws.Range("J114").Value = 926.75
$ws.Range("K114").Value = 933.5999999999999
$ws.Range("L114").Value = 2780.25
$ws.Range("M114").Value = 2320.4
$ws.Range("N114").Value = -9288.25

$ws.Range("H117").Value = 5839.095
$ws.Range("J117").Value = 6362.6313
$ws.Range("L117").Value = 19087.8939
$ws.Range("N117").Value = -25971.8939

$ws.Range("H118").Value = 2844
$ws.Range("I118").Value = 500
$ws.Range("J118").Value = 4016
$ws.Range("K118").Value = 1500
$ws.Range("L118").Value = 12048
$ws.Range("M118").Value = -257
$ws.Range("N118").Value = -14534

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34570.426
$ws.Range("I70").Value = 42620.652
$ws.Range("J70").Value = 4669.5713
$ws.Range("K70").Value = 42620.652
$ws.Range("L70").Value = 4669.5713
$ws.Range("M70").Value = -42350.652
$ws.Range("N70").Value = -5209.5713

$ws.Range("H73").Value = 34570.426
$ws.Range("I73").Value = 42620.652
$ws.Range("J73").Value = 4669.5713
$ws.Range("K73").Value = 42620.652
$ws.Range("L73").Value = 4669.5713
$ws.Range("M73").Value = -41684.652
$ws.Range("N73").Value = -6541.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 70007
$ws.Range("J23").Value = 70007
$ws.Range("L23").Value = 70007
$ws.Range("N23").Value = -70467

$ws.Range("H24").Value = 17933.334
$ws.Range("J24").Value = 17933.334
$ws.Range("L24").Value = 17933.334
$ws.Range("N24").Value = -18619.334

$ws.Range("H25").Value = 28000
$ws.Range("I25").Value = 2500
$ws.Range("J25").Value = 130000
$ws.Range("K25").Value = 2500
$ws.Range("L25").Value = 130000
$ws.Range("M25").Value = -2270
$ws.Range("N25").Value = -130460

$ws.Range("H35").Value = 5507.1665
$ws.Range("I35").Value = 810.75
$ws.Range("J35").Value = 14900
$ws.Range("K35").Value = 810.75
$ws.Range("L35").Value = 14900
$ws.Range("M35").Value = -474.75
$ws.Range("N35").Value = -15572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6774.7856
$ws.Range("J54").Value = 6752.077
$ws.Range("L54").Value = 6752.077
$ws.Range("N54").Value = -7792.077

$ws.Range("H81").Value = 122926.336
$ws.Range("I81").Value = 757.6
$ws.Range("K81").Value = 1515.2
$ws.Range("M81").Value = -454.2

$ws.Range("H84").Value = 122926.336
$ws.Range("I84").Value = 757.6
$ws.Range("K84").Value = 7576
$ws.Range("M84").Value = -2272

$ws.Range("H132").Value = 2792.2432
$ws.Range("I132").Value = 2794.0625
$ws.Range("K132").Value = 8382.1875
$ws.Range("M132").Value = -5852.1875
